# us_med_pnmcca_observatoire_comptage_terrain_meteo_2024-06-26.xlsx
# Commit: "Changed final format from csv to tsv to prevent any changes and
# problems while exchanging databases from a person to another. Other
# changes : corrected error from the meteo files since some errors were
# not corrected when the implementation of column content verification
# was done."
#
# Concretely this resolves to:
#   1) Fix two "Secteur" sector-name typos that had slipped through the
#      column-content verification pass:
#        "Torra di Murtella"  -> "A Torra di Murtella"  (rows 18-21)
#        "Maffalcu"           -> "Malfalcu"              (rows 30-33)
#   2) Widen column B ("Secteur") so the longer corrected labels are not
#      truncated on screen.
#   3) Leave the selection on the corrected "Malfalcu" block (B30:B33) as
#      a visual cue of the last edit made.
#   4) Shrink the saved window height slightly (leftover from resizing the
#      app window while reviewing the fix).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Correct the two sector-name typos -------------------------------
# "Torra di Murtella" appears in B18:B21
$ws.Range("B18").Value = "A Torra di Murtella"
$ws.Range("B19").Value = "A Torra di Murtella"
$ws.Range("B20").Value = "A Torra di Murtella"
$ws.Range("B21").Value = "A Torra di Murtella"

# "Maffalcu" appears in B30:B33
$ws.Range("B30").Value = "Malfalcu"
$ws.Range("B31").Value = "Malfalcu"
$ws.Range("B32").Value = "Malfalcu"
$ws.Range("B33").Value = "Malfalcu"

# --- 2) Widen column B to fit the corrected / longer labels -------------
$ws.Columns("B").ColumnWidth = 32.333333333333336

# --- 3) Put the selection on the last corrected block (B30:B33) ---------
$ws.Range("B30:B33").Select()

# --- 4) Shrink the recorded window height --------------------------------
$excel.ActiveWindow.Height = 609
